$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 12
$ws.Range("F5").Value = 50
$ws.Range("L5").Value = 30
$ws.Range("L14").Value = 9310
